$wb = $excel.ActiveWorkbook

# --- DatosCuenta: update QA smoke identifiers (row 2) ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQAJuneOne"
$wsCuenta.Range("B2").Value = "SmokeNameQAJuneOne"

# --- DatosMotor: remove the now-duplicate third data row ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A3:D3").Select()
$wsMotor.Rows("3").Delete()

# --- restore original active sheet / selection (DatosAP was active) ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Activate()
$wsAP.Range("F14").Select()
